# Insert a new weekly price record as the new row 207, pushing the existing
# rows 207:230 down to 208:231 (matching the "Fruta / hortaliza, semanal" update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 207. This shifts rows 207-230 down
# to 208-231 and carries over cell formatting/styles automatically.
$ws.Rows("207").Insert()

# Populate the newly inserted row 207 with the new data record.
$ws.Range("A207").Value2 = 1
$ws.Range("B207").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C207").Value2 = "Arica y Parinacota"
$ws.Range("D207").Value2 = 45106
$ws.Range("E207").Value2 = 15
$ws.Range("F207").Value2 = 100114001
$ws.Range("G207").Value2 = "Papa"
$ws.Range("H207").Value2 = "Cardinal"
$ws.Range("I207").Value2 = "1a nueva(o)"
$ws.Range("J207").Value2 = 1250
$ws.Range("K207").Value2 = 16000
$ws.Range("L207").Value2 = 17000
$ws.Range("M207").Value2 = 16520
$ws.Range("N207").Value2 = "`$/saco 25 kilos"
$ws.Range("O207").Value2 = "Región de Coquimbo"
$ws.Range("P207").Value2 = 661
$ws.Range("Q207").Value2 = 25
$ws.Range("R207").Value2 = "Hortaliza"
